# Update cryptocurrency price/volume snapshot data (coinranking.com scrape)
# Commit: "Updated cryptos list on Sat Mar 18 14:37:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.605.85"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "1.825.16"
$ws.Range("E3").Value = "  +5.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.67"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3539"
$ws.Range("E8").Value = "  +5.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.77"
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.240"
$ws.Range("E10").Value = "  +4.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07762"
$ws.Range("E11").Value = "  +4.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  +10.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.642"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").Value = "1.823.22"
$ws.Range("E15").Value = "  +5.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.209"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001128"
$ws.Range("E17").Value = "  +4.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06728"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  +4.94%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.70"
$ws.Range("E21").Value = "  +6.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.558"
$ws.Range("E22").Value = "  +6.53%  "
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").Value = "27.615.91"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.480"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.676"
$ws.Range("E26").Value = "  +11.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.13"
$ws.Range("E27").Value = "  +13.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.491"
$ws.Range("E28").Value = "  +6.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.06"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "2.027.29"
$ws.Range("E30").Value = "  +4.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.59"
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.364"
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.95"
$ws.Range("E34").Value = "  +8.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08831"
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.698"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.655"
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7051"
$ws.Range("E38").Value = "  +13.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.144"
$ws.Range("E39").Value = "  +7.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2268"
$ws.Range("E40").Value = "  +4.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06536"
$ws.Range("E41").Value = "  +4.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02407"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.303"
$ws.Range("E43").Value = "  +6.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.88"
$ws.Range("E44").Value = "  +5.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6632"
$ws.Range("E45").Value = "  +10.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9994"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.915"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.196"
$ws.Range("E48").Value = "  +7.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "133.73"
$ws.Range("E49").Value = "  +4.44%  "
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.62"
$ws.Range("E51").Value = "  +5.57%  "
